$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated statistics) ---
$ws.Cells.Item(3, 7).Value = 0.0086716268681249
$ws.Cells.Item(4, 7).Value = 0.0086716268681249
$ws.Cells.Item(5, 7).Value = 1130.93121191056
$ws.Cells.Item(6, 7).Value = 1130.93121191056
$ws.Cells.Item(7, 7).Value = 1130.93121191056
$ws.Cells.Item(8, 7).Value = 1130.93121191056
$ws.Cells.Item(9, 7).Value = 0.0114314858605896
$ws.Cells.Item(10, 7).Value = 0.0114314858605896
$ws.Cells.Item(13, 7).Value = 0.858054385964912
$ws.Cells.Item(13, 8).Value = 1.8002
$ws.Cells.Item(13, 9).Value = 1.59467
$ws.Cells.Item(13, 12).Value = 0.28875
$ws.Cells.Item(13, 13).Value = 1.3508
$ws.Cells.Item(14, 7).Value = 0.858054385964912
$ws.Cells.Item(14, 8).Value = 1.8002
$ws.Cells.Item(14, 9).Value = 1.59467
$ws.Cells.Item(14, 12).Value = 0.28875
$ws.Cells.Item(14, 13).Value = 1.3508
$ws.Cells.Item(17, 7).Value = 0.028856369323615
$ws.Cells.Item(18, 7).Value = 0.028856369323615
$ws.Cells.Item(20, 7).Value = 0.0088435378538622
$ws.Cells.Item(21, 7).Value = 0.0088435378538622
$ws.Cells.Item(22, 7).Value = 1653.73465131503
$ws.Cells.Item(22, 9).Value = 8713.03954
$ws.Cells.Item(23, 7).Value = 1653.73465131503
$ws.Cells.Item(23, 9).Value = 8713.03954
$ws.Cells.Item(24, 7).Value = 1653.73465131503
$ws.Cells.Item(24, 9).Value = 8713.03954
$ws.Cells.Item(25, 7).Value = 1653.73465131503
$ws.Cells.Item(25, 9).Value = 8713.03954
$ws.Cells.Item(26, 7).Value = 0.009392804095679301
$ws.Cells.Item(26, 12).Value = 0.00128
$ws.Cells.Item(27, 7).Value = 0.009392804095679301
$ws.Cells.Item(27, 12).Value = 0.00128
$ws.Cells.Item(30, 6).Value = 0.97085
$ws.Cells.Item(30, 7).Value = 0.903235
$ws.Cells.Item(30, 8).Value = 1.8002
$ws.Cells.Item(30, 9).Value = 1.6367
$ws.Cells.Item(30, 13).Value = 1.39455
$ws.Cells.Item(31, 6).Value = 0.97085
$ws.Cells.Item(31, 7).Value = 0.903235
$ws.Cells.Item(31, 8).Value = 1.8002
$ws.Cells.Item(31, 9).Value = 1.6367
$ws.Cells.Item(31, 13).Value = 1.39455
$ws.Cells.Item(37, 7).Value = 0.009277775271301599
$ws.Cells.Item(38, 7).Value = 0.009277775271301599
$ws.Cells.Item(39, 7).Value = 1328.63465131503
$ws.Cells.Item(40, 7).Value = 1328.63465131503
$ws.Cells.Item(41, 7).Value = 1328.63465131503
$ws.Cells.Item(42, 7).Value = 1328.63465131503
$ws.Cells.Item(43, 6).Value = 0.00185
$ws.Cells.Item(43, 7).Value = 0.0081078627018123
$ws.Cells.Item(43, 12).Value = 0.00128
$ws.Cells.Item(44, 6).Value = 0.00185
$ws.Cells.Item(44, 7).Value = 0.0081078627018123
$ws.Cells.Item(44, 12).Value = 0.00128
$ws.Cells.Item(47, 8).Value = 1.8002
$ws.Cells.Item(47, 9).Value = 1.70395
$ws.Cells.Item(47, 13).Value = 1.43685
$ws.Cells.Item(47, 14).Value = 1.62565
$ws.Cells.Item(48, 8).Value = 1.8002
$ws.Cells.Item(48, 9).Value = 1.70395
$ws.Cells.Item(48, 13).Value = 1.43685
$ws.Cells.Item(48, 14).Value = 1.62565
$ws.Cells.Item(54, 7).Value = 0.0100573307811641
$ws.Cells.Item(55, 7).Value = 0.0100573307811641
$ws.Cells.Item(56, 7).Value = 2146.76405218477
$ws.Cells.Item(56, 9).Value = 10982.24349
$ws.Cells.Item(57, 7).Value = 2146.76405218477
$ws.Cells.Item(57, 9).Value = 10982.24349
$ws.Cells.Item(58, 7).Value = 2146.76405218477
$ws.Cells.Item(58, 9).Value = 10982.24349
$ws.Cells.Item(59, 7).Value = 2146.76405218477
$ws.Cells.Item(59, 9).Value = 10982.24349
$ws.Cells.Item(60, 7).Value = 0.0098708586295222
$ws.Cells.Item(60, 12).Value = 0.00231
$ws.Cells.Item(61, 7).Value = 0.0098708586295222
$ws.Cells.Item(61, 12).Value = 0.00231
$ws.Cells.Item(64, 7).Value = 0.945301694915254
$ws.Cells.Item(64, 8).Value = 1.8002
$ws.Cells.Item(64, 9).Value = 1.70855
$ws.Cells.Item(64, 13).Value = 1.43762
$ws.Cells.Item(64, 14).Value = 1.62681
$ws.Cells.Item(65, 7).Value = 0.945301694915254
$ws.Cells.Item(65, 8).Value = 1.8002
$ws.Cells.Item(65, 9).Value = 1.70855
$ws.Cells.Item(65, 13).Value = 1.43762
$ws.Cells.Item(65, 14).Value = 1.62681
$ws.Cells.Item(71, 7).Value = 0.0102098731540454
$ws.Cells.Item(72, 7).Value = 0.0102098731540454
$ws.Cells.Item(73, 7).Value = 2154.25998411899
$ws.Cells.Item(73, 9).Value = 11262.44745
$ws.Cells.Item(74, 7).Value = 2154.25998411899
$ws.Cells.Item(74, 9).Value = 11262.44745
$ws.Cells.Item(75, 7).Value = 2154.25998411899
$ws.Cells.Item(75, 9).Value = 11262.44745
$ws.Cells.Item(76, 7).Value = 2154.25998411899
$ws.Cells.Item(76, 9).Value = 11262.44745
$ws.Cells.Item(77, 7).Value = 0.0084761342969707
$ws.Cells.Item(77, 12).Value = 0.00351
$ws.Cells.Item(78, 7).Value = 0.0084761342969707
$ws.Cells.Item(78, 12).Value = 0.00351
$ws.Cells.Item(81, 6).Value = 0.9757
$ws.Cells.Item(81, 9).Value = 1.69205
$ws.Cells.Item(81, 14).Value = 1.62681
$ws.Cells.Item(82, 6).Value = 0.9757
$ws.Cells.Item(82, 9).Value = 1.69205
$ws.Cells.Item(82, 14).Value = 1.62681
$ws.Cells.Item(88, 7).Value = 0.0099623461832844
$ws.Cells.Item(89, 7).Value = 0.0099623461832844
$ws.Cells.Item(94, 7).Value = 0.009206912688541
$ws.Cells.Item(95, 7).Value = 0.009206912688541
$ws.Cells.Item(105, 7).Value = 0.0102579472878428
$ws.Cells.Item(106, 7).Value = 0.0102579472878428
$ws.Cells.Item(111, 7).Value = 0.0106984248320634
$ws.Cells.Item(112, 7).Value = 0.0106984248320634
$ws.Cells.Item(122, 7).Value = 0.0113096714257738
$ws.Cells.Item(123, 7).Value = 0.0113096714257738
$ws.Cells.Item(128, 7).Value = 0.0128148457219132
$ws.Cells.Item(129, 7).Value = 0.0128148457219132

# --- Append new rows 138-154 (2019-2023 period results) ---
# Row 138
$ws.Cells.Item(138, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(138, 2).Value = "Visual Clarity (Sediment class 4)"
$ws.Cells.Item(138, 3).Value = "A"
$ws.Cells.Item(138, 4).Value = "2019 - 2023"
$ws.Cells.Item(138, 5).Value = "RepSite"
$ws.Cells.Item(138, 6).Value = 1.75
$ws.Cells.Item(138, 7).Value = 1.5790243902439
$ws.Cells.Item(138, 8).Value = 4.1
$ws.Cells.Item(138, 9).Value = 3.645
$ws.Cells.Item(138, 12).Value = 1.8
$ws.Cells.Item(138, 13).Value = 2.8201
$ws.Cells.Item(138, 14).Value = 3.327
$ws.Cells.Item(138, 15).Value = 1827830.498
$ws.Cells.Item(138, 16).Value = 5496951.677
$ws.Cells.Item(138, 17).Value = "Tararua District"
$ws.Cells.Item(138, 18).Value = "Manawatū"
$ws.Cells.Item(138, 19).Value = "Mangatainoka"
$ws.Cells.Item(138, 20).Value = "Mana_8d"
$ws.Cells.Item(138, 21).Value = "m"
# Row 139
$ws.Cells.Item(139, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(139, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(139, 3).Value = "B"
$ws.Cells.Item(139, 4).Value = "2019 - 2023"
$ws.Cells.Item(139, 5).Value = "RepSite"
$ws.Cells.Item(139, 6).Value = 0.011
$ws.Cells.Item(139, 7).Value = 0.0117620498761844
$ws.Cells.Item(139, 8).Value = 0.037
$ws.Cells.Item(139, 9).Value = 0.0262
$ws.Cells.Item(139, 12).Value = 0.0105
$ws.Cells.Item(139, 13).Value = 0.01647
$ws.Cells.Item(139, 14).Value = 0.02156
$ws.Cells.Item(139, 15).Value = 1827830.498
$ws.Cells.Item(139, 16).Value = 5496951.677
$ws.Cells.Item(139, 17).Value = "Tararua District"
$ws.Cells.Item(139, 18).Value = "Manawatū"
$ws.Cells.Item(139, 19).Value = "Mangatainoka"
$ws.Cells.Item(139, 20).Value = "Mana_8d"
$ws.Cells.Item(139, 21).Value = "mg/L"
# Row 140
$ws.Cells.Item(140, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(140, 2).Value = "DRP (Median)"
$ws.Cells.Item(140, 3).Value = "C"
$ws.Cells.Item(140, 4).Value = "2019 - 2023"
$ws.Cells.Item(140, 5).Value = "RepSite"
$ws.Cells.Item(140, 6).Value = 0.011
$ws.Cells.Item(140, 7).Value = 0.0117620498761844
$ws.Cells.Item(140, 8).Value = 0.037
$ws.Cells.Item(140, 9).Value = 0.0262
$ws.Cells.Item(140, 12).Value = 0.0105
$ws.Cells.Item(140, 13).Value = 0.01647
$ws.Cells.Item(140, 14).Value = 0.02156
$ws.Cells.Item(140, 15).Value = 1827830.498
$ws.Cells.Item(140, 16).Value = 5496951.677
$ws.Cells.Item(140, 17).Value = "Tararua District"
$ws.Cells.Item(140, 18).Value = "Manawatū"
$ws.Cells.Item(140, 19).Value = "Mangatainoka"
$ws.Cells.Item(140, 20).Value = "Mana_8d"
$ws.Cells.Item(140, 21).Value = "mg/L"
# Row 141
$ws.Cells.Item(141, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(141, 2).Value = "E coli (>260)"
$ws.Cells.Item(141, 3).Value = "E"
$ws.Cells.Item(141, 4).Value = "2019 - 2023"
$ws.Cells.Item(141, 5).Value = "RepSite"
$ws.Cells.Item(141, 6).Value = 466
$ws.Cells.Item(141, 7).Value = 5464.44827586207
$ws.Cells.Item(141, 8).Value = 110000
$ws.Cells.Item(141, 9).Value = 32324
$ws.Cells.Item(141, 10).Value = 43.1034482758621
$ws.Cells.Item(141, 11).Value = 62.0689655172414
$ws.Cells.Item(141, 12).Value = 280
$ws.Cells.Item(141, 13).Value = 2850.04
$ws.Cells.Item(141, 14).Value = 27160
$ws.Cells.Item(141, 15).Value = 1827830.498
$ws.Cells.Item(141, 16).Value = 5496951.677
$ws.Cells.Item(141, 17).Value = "Tararua District"
$ws.Cells.Item(141, 18).Value = "Manawatū"
$ws.Cells.Item(141, 19).Value = "Mangatainoka"
$ws.Cells.Item(141, 20).Value = "Mana_8d"
$ws.Cells.Item(141, 21).Value = "% exceedances over 260/100 mL"
# Row 142
$ws.Cells.Item(142, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(142, 2).Value = "E coli (>540)"
$ws.Cells.Item(142, 3).Value = "E"
$ws.Cells.Item(142, 4).Value = "2019 - 2023"
$ws.Cells.Item(142, 5).Value = "RepSite"
$ws.Cells.Item(142, 6).Value = 466
$ws.Cells.Item(142, 7).Value = 5464.44827586207
$ws.Cells.Item(142, 8).Value = 110000
$ws.Cells.Item(142, 9).Value = 32324
$ws.Cells.Item(142, 10).Value = 43.1034482758621
$ws.Cells.Item(142, 11).Value = 62.0689655172414
$ws.Cells.Item(142, 12).Value = 280
$ws.Cells.Item(142, 13).Value = 2850.04
$ws.Cells.Item(142, 14).Value = 27160
$ws.Cells.Item(142, 15).Value = 1827830.498
$ws.Cells.Item(142, 16).Value = 5496951.677
$ws.Cells.Item(142, 17).Value = "Tararua District"
$ws.Cells.Item(142, 18).Value = "Manawatū"
$ws.Cells.Item(142, 19).Value = "Mangatainoka"
$ws.Cells.Item(142, 20).Value = "Mana_8d"
$ws.Cells.Item(142, 21).Value = "% exceedances over 540/100 mL"
# Row 143
$ws.Cells.Item(143, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(143, 2).Value = "E coli (Median)"
$ws.Cells.Item(143, 3).Value = "E"
$ws.Cells.Item(143, 4).Value = "2019 - 2023"
$ws.Cells.Item(143, 5).Value = "RepSite"
$ws.Cells.Item(143, 6).Value = 466
$ws.Cells.Item(143, 7).Value = 5464.44827586207
$ws.Cells.Item(143, 8).Value = 110000
$ws.Cells.Item(143, 9).Value = 32324
$ws.Cells.Item(143, 10).Value = 43.1034482758621
$ws.Cells.Item(143, 11).Value = 62.0689655172414
$ws.Cells.Item(143, 12).Value = 280
$ws.Cells.Item(143, 13).Value = 2850.04
$ws.Cells.Item(143, 14).Value = 27160
$ws.Cells.Item(143, 15).Value = 1827830.498
$ws.Cells.Item(143, 16).Value = 5496951.677
$ws.Cells.Item(143, 17).Value = "Tararua District"
$ws.Cells.Item(143, 18).Value = "Manawatū"
$ws.Cells.Item(143, 19).Value = "Mangatainoka"
$ws.Cells.Item(143, 20).Value = "Mana_8d"
$ws.Cells.Item(143, 21).Value = "E. coli/100 mL"
# Row 144
$ws.Cells.Item(144, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(144, 2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(144, 3).Value = "E"
$ws.Cells.Item(144, 4).Value = "2019 - 2023"
$ws.Cells.Item(144, 5).Value = "RepSite"
$ws.Cells.Item(144, 6).Value = 466
$ws.Cells.Item(144, 7).Value = 5464.44827586207
$ws.Cells.Item(144, 8).Value = 110000
$ws.Cells.Item(144, 9).Value = 32324
$ws.Cells.Item(144, 10).Value = 43.1034482758621
$ws.Cells.Item(144, 11).Value = 62.0689655172414
$ws.Cells.Item(144, 12).Value = 280
$ws.Cells.Item(144, 13).Value = 2850.04
$ws.Cells.Item(144, 14).Value = 27160
$ws.Cells.Item(144, 15).Value = 1827830.498
$ws.Cells.Item(144, 16).Value = 5496951.677
$ws.Cells.Item(144, 17).Value = "Tararua District"
$ws.Cells.Item(144, 18).Value = "Manawatū"
$ws.Cells.Item(144, 19).Value = "Mangatainoka"
$ws.Cells.Item(144, 20).Value = "Mana_8d"
$ws.Cells.Item(144, 21).Value = "E. coli/100 mL"
# Row 145
$ws.Cells.Item(145, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(145, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(145, 3).Value = "A"
$ws.Cells.Item(145, 4).Value = "2019 - 2023"
$ws.Cells.Item(145, 5).Value = "RepSite"
$ws.Cells.Item(145, 6).Value = 0.0074
$ws.Cells.Item(145, 7).Value = 0.0134700020211978
$ws.Cells.Item(145, 8).Value = 0.052395069180857
$ws.Cells.Item(145, 9).Value = 0.04535
$ws.Cells.Item(145, 12).Value = 0.00903
$ws.Cells.Item(145, 13).Value = 0.02541
$ws.Cells.Item(145, 14).Value = 0.0428
$ws.Cells.Item(145, 15).Value = 1827830.498
$ws.Cells.Item(145, 16).Value = 5496951.677
$ws.Cells.Item(145, 17).Value = "Tararua District"
$ws.Cells.Item(145, 18).Value = "Manawatū"
$ws.Cells.Item(145, 19).Value = "Mangatainoka"
$ws.Cells.Item(145, 20).Value = "Mana_8d"
$ws.Cells.Item(145, 21).Value = "mg NH4-N/L"
# Row 146
$ws.Cells.Item(146, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(146, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(146, 3).Value = "A"
$ws.Cells.Item(146, 4).Value = "2019 - 2023"
$ws.Cells.Item(146, 5).Value = "RepSite"
$ws.Cells.Item(146, 6).Value = 0.0074
$ws.Cells.Item(146, 7).Value = 0.0134700020211978
$ws.Cells.Item(146, 8).Value = 0.052395069180857
$ws.Cells.Item(146, 9).Value = 0.04535
$ws.Cells.Item(146, 12).Value = 0.00903
$ws.Cells.Item(146, 13).Value = 0.02541
$ws.Cells.Item(146, 14).Value = 0.0428
$ws.Cells.Item(146, 15).Value = 1827830.498
$ws.Cells.Item(146, 16).Value = 5496951.677
$ws.Cells.Item(146, 17).Value = "Tararua District"
$ws.Cells.Item(146, 18).Value = "Manawatū"
$ws.Cells.Item(146, 19).Value = "Mangatainoka"
$ws.Cells.Item(146, 20).Value = "Mana_8d"
$ws.Cells.Item(146, 21).Value = "mg NH4-N/L"
# Row 147
$ws.Cells.Item(147, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(147, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(147, 3).Value = "B"
$ws.Cells.Item(147, 4).Value = "2019 - 2023"
$ws.Cells.Item(147, 5).Value = "RepSite"
$ws.Cells.Item(147, 6).Value = 0.759
$ws.Cells.Item(147, 7).Value = 0.819457627118644
$ws.Cells.Item(147, 8).Value = 1.86
$ws.Cells.Item(147, 9).Value = 1.6325
$ws.Cells.Item(147, 12).Value = 0.5125
$ws.Cells.Item(147, 13).Value = 1.2076
$ws.Cells.Item(147, 14).Value = 1.4268
$ws.Cells.Item(147, 15).Value = 1827830.498
$ws.Cells.Item(147, 16).Value = 5496951.677
$ws.Cells.Item(147, 17).Value = "Tararua District"
$ws.Cells.Item(147, 18).Value = "Manawatū"
$ws.Cells.Item(147, 19).Value = "Mangatainoka"
$ws.Cells.Item(147, 20).Value = "Mana_8d"
$ws.Cells.Item(147, 21).Value = "mg NO3-N/L"
# Row 148
$ws.Cells.Item(148, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(148, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(148, 3).Value = "A"
$ws.Cells.Item(148, 4).Value = "2019 - 2023"
$ws.Cells.Item(148, 5).Value = "RepSite"
$ws.Cells.Item(148, 6).Value = 0.759
$ws.Cells.Item(148, 7).Value = 0.819457627118644
$ws.Cells.Item(148, 8).Value = 1.86
$ws.Cells.Item(148, 9).Value = 1.6325
$ws.Cells.Item(148, 12).Value = 0.5125
$ws.Cells.Item(148, 13).Value = 1.2076
$ws.Cells.Item(148, 14).Value = 1.4268
$ws.Cells.Item(148, 15).Value = 1827830.498
$ws.Cells.Item(148, 16).Value = 5496951.677
$ws.Cells.Item(148, 17).Value = "Tararua District"
$ws.Cells.Item(148, 18).Value = "Manawatū"
$ws.Cells.Item(148, 19).Value = "Mangatainoka"
$ws.Cells.Item(148, 20).Value = "Mana_8d"
$ws.Cells.Item(148, 21).Value = "mg NO3-N/L"
# Row 149
$ws.Cells.Item(149, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(149, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(149, 4).Value = "2019 - 2023"
$ws.Cells.Item(149, 5).Value = "RepSite"
$ws.Cells.Item(149, 6).Value = 0.79
$ws.Cells.Item(149, 7).Value = 0.848983050847458
$ws.Cells.Item(149, 8).Value = 1.89
$ws.Cells.Item(149, 9).Value = 1.657
$ws.Cells.Item(149, 12).Value = 0.545
$ws.Cells.Item(149, 13).Value = 1.237
$ws.Cells.Item(149, 14).Value = 1.517
$ws.Cells.Item(149, 15).Value = 1827830.498
$ws.Cells.Item(149, 16).Value = 5496951.677
$ws.Cells.Item(149, 17).Value = "Tararua District"
$ws.Cells.Item(149, 18).Value = "Manawatū"
$ws.Cells.Item(149, 19).Value = "Mangatainoka"
$ws.Cells.Item(149, 20).Value = "Mana_8d"
$ws.Cells.Item(149, 21).Value = "g/m3"
# Row 150
$ws.Cells.Item(150, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(150, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(150, 4).Value = "2019 - 2023"
$ws.Cells.Item(150, 5).Value = "RepSite"
$ws.Cells.Item(150, 6).Value = 0.79
$ws.Cells.Item(150, 7).Value = 0.848983050847458
$ws.Cells.Item(150, 8).Value = 1.89
$ws.Cells.Item(150, 9).Value = 1.657
$ws.Cells.Item(150, 12).Value = 0.545
$ws.Cells.Item(150, 13).Value = 1.237
$ws.Cells.Item(150, 14).Value = 1.517
$ws.Cells.Item(150, 15).Value = 1827830.498
$ws.Cells.Item(150, 16).Value = 5496951.677
$ws.Cells.Item(150, 17).Value = "Tararua District"
$ws.Cells.Item(150, 18).Value = "Manawatū"
$ws.Cells.Item(150, 19).Value = "Mangatainoka"
$ws.Cells.Item(150, 20).Value = "Mana_8d"
$ws.Cells.Item(150, 21).Value = "g/m3"
# Row 151
$ws.Cells.Item(151, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(151, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(151, 4).Value = "2019 - 2023"
$ws.Cells.Item(151, 5).Value = "RepSite"
$ws.Cells.Item(151, 6).Value = 1.13
$ws.Cells.Item(151, 7).Value = 1.15474576271186
$ws.Cells.Item(151, 8).Value = 2.08
$ws.Cells.Item(151, 9).Value = 1.9015
$ws.Cells.Item(151, 12).Value = 0.87
$ws.Cells.Item(151, 13).Value = 1.5741
$ws.Cells.Item(151, 14).Value = 1.7568
$ws.Cells.Item(151, 15).Value = 1827830.498
$ws.Cells.Item(151, 16).Value = 5496951.677
$ws.Cells.Item(151, 17).Value = "Tararua District"
$ws.Cells.Item(151, 18).Value = "Manawatū"
$ws.Cells.Item(151, 19).Value = "Mangatainoka"
$ws.Cells.Item(151, 20).Value = "Mana_8d"
$ws.Cells.Item(151, 21).Value = "g/m3"
# Row 152
$ws.Cells.Item(152, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(152, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(152, 4).Value = "2019 - 2023"
$ws.Cells.Item(152, 5).Value = "RepSite"
$ws.Cells.Item(152, 6).Value = 1.13
$ws.Cells.Item(152, 7).Value = 1.15474576271186
$ws.Cells.Item(152, 8).Value = 2.08
$ws.Cells.Item(152, 9).Value = 1.9015
$ws.Cells.Item(152, 12).Value = 0.87
$ws.Cells.Item(152, 13).Value = 1.5741
$ws.Cells.Item(152, 14).Value = 1.7568
$ws.Cells.Item(152, 15).Value = 1827830.498
$ws.Cells.Item(152, 16).Value = 5496951.677
$ws.Cells.Item(152, 17).Value = "Tararua District"
$ws.Cells.Item(152, 18).Value = "Manawatū"
$ws.Cells.Item(152, 19).Value = "Mangatainoka"
$ws.Cells.Item(152, 20).Value = "Mana_8d"
$ws.Cells.Item(152, 21).Value = "g/m3"
# Row 153
$ws.Cells.Item(153, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(153, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(153, 4).Value = "2019 - 2023"
$ws.Cells.Item(153, 5).Value = "RepSite"
$ws.Cells.Item(153, 6).Value = 0.031
$ws.Cells.Item(153, 7).Value = 0.0570508474576271
$ws.Cells.Item(153, 8).Value = 0.336
$ws.Cells.Item(153, 9).Value = 0.21615
$ws.Cells.Item(153, 12).Value = 0.0295
$ws.Cells.Item(153, 13).Value = 0.08863
$ws.Cells.Item(153, 14).Value = 0.17788
$ws.Cells.Item(153, 15).Value = 1827830.498
$ws.Cells.Item(153, 16).Value = 5496951.677
$ws.Cells.Item(153, 17).Value = "Tararua District"
$ws.Cells.Item(153, 18).Value = "Manawatū"
$ws.Cells.Item(153, 19).Value = "Mangatainoka"
$ws.Cells.Item(153, 20).Value = "Mana_8d"
$ws.Cells.Item(153, 21).Value = "g/m3"
# Row 154
$ws.Cells.Item(154, 1).Value = "Ngatahaka Stream at u/s Makakahi Confl"
$ws.Cells.Item(154, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(154, 4).Value = "2019 - 2023"
$ws.Cells.Item(154, 5).Value = "RepSite"
$ws.Cells.Item(154, 6).Value = 0.031
$ws.Cells.Item(154, 7).Value = 0.0570508474576271
$ws.Cells.Item(154, 8).Value = 0.336
$ws.Cells.Item(154, 9).Value = 0.21615
$ws.Cells.Item(154, 12).Value = 0.0295
$ws.Cells.Item(154, 13).Value = 0.08863
$ws.Cells.Item(154, 14).Value = 0.17788
$ws.Cells.Item(154, 15).Value = 1827830.498
$ws.Cells.Item(154, 16).Value = 5496951.677
$ws.Cells.Item(154, 17).Value = "Tararua District"
$ws.Cells.Item(154, 18).Value = "Manawatū"
$ws.Cells.Item(154, 19).Value = "Mangatainoka"
$ws.Cells.Item(154, 20).Value = "Mana_8d"
$ws.Cells.Item(154, 21).Value = "g/m3"
